$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'254.24"
$ws.Range("E2").Value = "'3.25%"
$ws.Range("D3").Value = "'27.99"
$ws.Range("E3").Value = "'-5.57%"
$ws.Range("D4").Value = "'5.321"
$ws.Range("E4").Value = "'3.12%"
$ws.Range("D5").Value = "'0.05839"
$ws.Range("E5").Value = "'0.61%"
$ws.Range("D6").Value = "'6.706"
$ws.Range("E6").Value = "'0.72%"
$ws.Range("D7").Value = "'0.8638"
$ws.Range("E7").Value = "'1.48%"
$ws.Range("D8").Value = "'0.9119"
$ws.Range("E8").Value = "'5.36%"
$ws.Range("D9").Value = "'0.1428"
$ws.Range("E9").Value = "'3.46%"
$ws.Range("D10").Value = "'0.07173"
$ws.Range("E10").Value = "'1.19%"
$ws.Range("D11").Value = "'0.03210"
$ws.Range("E11").Value = "'0.36%"
$ws.Range("D12").Value = "'0.09224"
$ws.Range("E12").Value = "'-1.61%"
$ws.Range("D13").Value = "'0.001553"
$ws.Range("E13").Value = "'2.03%"
$ws.Range("D14").Value = "'0.0006039"
$ws.Range("E14").Value = "'0.88%"
$ws.Range("D15").Value = "'0.006063"
$ws.Range("E15").Value = "'-1.45%"
$ws.Range("D16").Value = "'3.498"
$ws.Range("E16").Value = "'0.33%"
$ws.Range("E17").Value = "'0.84%"
$ws.Range("E18").Value = "'0.12%"
$ws.Range("E19").Value = "'-0.89%"
$ws.Range("D20").Value = "'0.03454"
$ws.Range("E20").Value = "'2.62%"
$ws.Range("D21").Value = "'0.1334"
$ws.Range("E21").Value = "'4.06%"
$ws.Range("D22").Value = "'3.526"
$ws.Range("E22").Value = "'1.38%"
$ws.Range("D23").Value = "'0.04147"
$ws.Range("E23").Value = "'0.19%"
$ws.Range("D24").Value = "'0.1379"
$ws.Range("E24").Value = "'-0.10%"
$ws.Range("D26").Value = "'0.001225"
$ws.Range("E26").Value = "'-0.15%"
$ws.Range("D27").Value = "'0.0001199"
$ws.Range("E27").Value = "'-0.77%"
$ws.Range("D28").Value = "'0.0001938"
$ws.Range("E28").Value = "'34.23%"
$ws.Range("D40").Value = "'0.03858"
$ws.Range("E40").Value = "'3.06%"
$ws.Range("D41").Value = "'0.1098"
$ws.Range("E41").Value = "'2.57%"
$ws.Range("D42").Value = "'0.002199"
$ws.Range("E42").Value = "'0.05%"
$ws.Range("D43").Value = "'0.002948"
$ws.Range("E43").Value = "'-48.58%"
$ws.Range("D44").Value = "'0.01088"
$ws.Range("E44").Value = "'18.46%"
$ws.Range("D45").Value = "'0.00005240"
$ws.Range("E45").Value = "'-1.13%"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("D47").Value = "'0.08984"
$ws.Range("E47").Value = "'55.07%"
$ws.Range("E48").Value = "'-0.89%"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E50").Value = "'0.05%"
